# Updates cryptos list prices (col D) and 1h volume % (col E) for Sheet1.
# Price cells that look like plain decimals (e.g. "375.93") are prefixed
# with a leading apostrophe so Excel stores them as literal text instead
# of auto-converting to a number (matches the source data's text cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.592.04"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.936.23"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("D5").Value = "'375.93"
$ws.Range("E5").Value = "  +6.25%  "
$ws.Range("D6").Value = "'104.22"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("E7").Value = "  -2.88%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  -3.98%  "
$ws.Range("D10").Value = "'37.03"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "'0.0839"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").Value = "'18.38"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "3.399.05"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "'7.39"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "2.936.84"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "'0.933"
$ws.Range("E17").Value = "  -8.09%  "
$ws.Range("D18").Value = "51.536.84"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "'3.43"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "'13.00"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").Value = "'68.39"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "'262.41"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").Value = "'4.14"
$ws.Range("E27").Value = "  -4.88%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").Value = "'25.84"
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("D31").Value = "'7.02"
$ws.Range("E31").Value = "  +8.54%  "
$ws.Range("D32").Value = "'0.102"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("D33").Value = "'9.84"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("D35").Value = "'51.09"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'34.20"
$ws.Range("E36").Value = "  -4.90%  "
$ws.Range("D37").Value = "'0.0430"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  -8.66%  "
$ws.Range("D40").Value = "'17.04"
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "'2.61"
$ws.Range("E41").Value = "  -7.53%  "
$ws.Range("D42").Value = "'1.83"
$ws.Range("E42").Value = "  -6.32%  "
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'122.04"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("E46").Value = "  -5.49%  "
$ws.Range("E47").Value = "  +11.85%  "
$ws.Range("D48").Value = "2.027.51"
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  -4.49%  "
$ws.Range("D51").Value = "3.215.59"
$ws.Range("E51").Value = "  -2.68%  "
